# Apply the "Added updated user and fixture list" commit:
#  - Round of 16 results (rows 38-45): fill in Home_Score (J) / Away_Score (K)
#  - Two new Quarter-final fixture days appended (rows 46-49)
#  - Workbook-level defined names (UNI_* helper names) removed
#  - Selection moved on to the newly added row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Remove the UNI_* hidden defined names from the workbook.
# ---------------------------------------------------------------------------
$nameCount = $wb.Names.Count
for ($i = 0; $i -lt $nameCount; $i++) {
    $wb.Names.Item(1).Delete()
}

# ---------------------------------------------------------------------------
# 2) Fill in the Round of 16 results for the already-listed fixtures.
# ---------------------------------------------------------------------------
$r16Results = @{
    38 = @(2, 0)   # Switzerland 2-0 Italy
    39 = @(2, 0)   # Germany 2-0 Denmark
    40 = @(1, 1)   # England 1-1 Slovakia
    41 = @(4, 1)   # Spain 4-1 Georgia
    42 = @(1, 0)   # France 1-0 Belgium
    43 = @(0, 0)   # Portugal 0-0 Slovenia
    44 = @(0, 3)   # Romania 0-3 Netherlands
    45 = @(1, 2)   # Austria 1-2 Turkey
}

foreach ($row in $r16Results.Keys) {
    $scores = $r16Results[$row]
    $ws.Cells.Item($row, 10).Value = $scores[0]   # column J - Home_Score
    $ws.Cells.Item($row, 11).Value = $scores[1]   # column K - Away_Score
}

# ---------------------------------------------------------------------------
# 3) Append the Quarter-final fixtures (rows 46-49).
# ---------------------------------------------------------------------------
$newFixtures = @(
    @{ Row = 46; Day = "Fri"; Date = "Jul 05, 2024"; Time = "18:00:00"; Home = "Spain";       Away = "Germany";     Venue = "Stuttgart" }
    @{ Row = 47; Day = "Fri"; Date = "Jul 05, 2024"; Time = "21:00:00"; Home = "Portugal";     Away = "France";      Venue = "Hamburg" }
    @{ Row = 48; Day = "Sat"; Date = "Jul 06, 2024"; Time = "18:00:00"; Home = "England";      Away = "Switzerland"; Venue = "Düsseldorf" }
    @{ Row = 49; Day = "Sat"; Date = "Jul 06, 2024"; Time = "21:00:00"; Home = "Netherlands";  Away = "Turkey";      Venue = "Berlin" }
)

foreach ($fixture in $newFixtures) {
    $row = $fixture.Row
    $ws.Cells.Item($row, 1).Value = $fixture.Day     # A - Day
    $ws.Cells.Item($row, 2).Value = $fixture.Date    # B - Date
    $ws.Cells.Item($row, 3).Value = $fixture.Time    # C - Time
    $ws.Cells.Item($row, 4).Value = $fixture.Home    # D - Home
    $ws.Cells.Item($row, 7).Value = $fixture.Away    # G - Away
    $ws.Cells.Item($row, 8).Value = $fixture.Venue   # H - Venue
}

# ---------------------------------------------------------------------------
# 4) Move the active selection onto the newly-added data, mirroring the
#    author's final cursor position in the saved workbook.
# ---------------------------------------------------------------------------
$null = $ws.Range("A16").Select()
$null = $ws.Range("J48").Select()
